$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New test "Check_ThatUserCanSearchByAnAlreadyExist_GF_Name" renames the
# GF_synonyms column header (H1) to "gfsynonyms" and keeps the data row's
# value (H2) as "GF_synonyms".
$ws.Range("H1").Value = "gfsynonyms"
$ws.Range("H2").Value = "GF_synonyms"

# Reflect the updated selection/view used while validating the change.
$ws.Range("C16").Select()
